# gsl2008, ssy19, phd2014 MethodReport_en
#
# "attachments" sheet: insert a new row for the English method report
# ("gsl2008_MethodReport_en.pdf"), mirroring the existing German method
# report row (row 2), and leave "attachments" as the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attachments")

# Insert a new row above the current row 3 ("...Overview_en.pdf"), pushing
# it - and the blank formatting row after it - down by one row. Excel/COM
# inherits the row-2 formatting for the freshly inserted row automatically.
$ws.Rows.Item(3).Insert() | Out-Null

# Populate the new row with the English method-report metadata, parallel to
# the existing German entry on row 2.
$ws.Range("A3").Value = "gsl2008_MethodReport_en.pdf"
$ws.Range("B3").Value = "Daten- und Methodenbericht"
$ws.Range("C3").Value = "Method Report"
$ws.Range("D3").Value = "DZHW Panel Study of School Leavers with a Higher Education Entrance Qualification 2008"
$ws.Range("E3").Value = "Daten- und Methodenbericht zu den Erhebungen des Studienberechtigtenjahrgangs 2008 (1. bis 3. Befragungswelle). Version 1.0.0"
$ws.Range("F3").Value = "Data and methods report for the surveys on School Leavers with a Higher Education Entrance Qualification of 2008. Version 1.0.0"
$ws.Range("G3").Value = "en"

# Make "attachments" the active/selected sheet (tab), with the blank
# formatting cell below the table (now on row 7) as the active selection.
$ws.Activate() | Out-Null
$ws.Range("D7").Select() | Out-Null
